# Insert a new weekly price record as row 34 ("Fruta / hortaliza, semanal").
# Every row below shifts down by one (old row 34 -> new row 35, ... old row 91 -> new row 92),
# and the new row 34 carries a fresh "Poroto verde" quote for Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 34, pushing rows 34:91 down to 35:92.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Macroferia Regional de Talca"
$ws.Range("C34").Value = "Maule"
$ws.Range("D34").Value = 44469
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = 100112031
$ws.Range("G34").Value = "Poroto verde"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 30000
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = 30000
$ws.Range("N34").Value = "$/malla 25 kilos"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 1200
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
